$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.132.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.44%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.842.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.41%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.82%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6870'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.96%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3015'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.81%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07449'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.63%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07650'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.26%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.839.20'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.55%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.057'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.57%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6824'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.73%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '87.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.99%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.173'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.79%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.122.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.35%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008143'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.17%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.077.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.81%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.27'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.76%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.76%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.387'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.95%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9994'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.02%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1451'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.22%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.47%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.759'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.94%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.514'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.74%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.279'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.13%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.135'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.22%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.192'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.96%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05231'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7579'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.58%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.848'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.99%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.133'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.36%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.685'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.35%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.301.05'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01834'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.12%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.723'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9299'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.55%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.860'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.55%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '104.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.57%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9992'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '

$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.981.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.53%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.09%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5196'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.20%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.480'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.54%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.767'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.13%  '

$ws.Range('B50').Value = 'XinFinNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07374'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +17.16%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05942'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.63%  '
